$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new September transaction-log entry was recorded at the top of the list.
# Insert a new row above row 47 so every existing entry (rows 47-193, plus
# the trailing "Broadband" label row at 193) shifts down by one -- this is
# exactly what the diff shows: each row's September Details/Date pair now
# holds the value that used to belong to the row above it, and the sheet's
# used range grows from A1:Y193 to A1:Y194.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row with the new entry's September details.
$ws.Range("R47").Value = "dispute"
$ws.Range("S47").Value = "2024-09-23 07:05:20"
